$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 171.82353
$ws.Range("J42").Value = 360.25
$ws.Range("L42").Value = 1080.75
$ws.Range("N42").Value = -1540.75
$ws.Range("H43").Value = 12629.4
$ws.Range("I43").Value = 12629.4
$ws.Range("K43").Value = 12629.4
$ws.Range("M43").Value = -12560.4
$ws.Range("H62").Value = 8336155
$ws.Range("I62").Value = 10418944
$ws.Range("J62").Value = 4999.5
$ws.Range("K62").Value = 10418944
$ws.Range("L62").Value = 4999.5
$ws.Range("M62").Value = -10418320
$ws.Range("N62").Value = -6247.5
$ws.Range("H65").Value = 8336155
$ws.Range("I65").Value = 10418944
$ws.Range("J65").Value = 4999.5
$ws.Range("K65").Value = 52094720
$ws.Range("L65").Value = 24997.5
$ws.Range("M65").Value = -52091600
$ws.Range("N65").Value = -31237.5
$ws.Range("H98").Value = 3119
$ws.Range("I98").Value = 3119
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 3119
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -1621
$ws.Range("N98").ClearContents()
$ws.Range("H106").Value = 3169.8948
$ws.Range("I106").Value = 3035.2
$ws.Range("J106").Value = 3675
$ws.Range("K106").Value = 3035.2
$ws.Range("L106").Value = 3675
$ws.Range("M106").Value = -2404.2
$ws.Range("N106").Value = -4937
$ws.Range("H122").Value = 3119
$ws.Range("I122").Value = 3119
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9357
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6907
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 1633.2
$ws.Range("I132").Value = 1279.7
$ws.Range("K132").Value = 3839.1
$ws.Range("M132").Value = -1309.1
$ws.Range("H138").Value = 3555.889
$ws.Range("I138").Value = 1857
$ws.Range("J138").Value = 5679.5
$ws.Range("K138").Value = 5571
$ws.Range("L138").Value = 17038.5
$ws.Range("M138").Value = -431
$ws.Range("N138").Value = -27318.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4330777.5
$ws.Range("I2").Value = 4786433
$ws.Range("J2").Value = 2050
$ws.Range("K2").Value = 4786433
$ws.Range("L2").Value = 2050
$ws.Range("M2").Value = -4786320
$ws.Range("N2").Value = -2276
$ws.Range("H45").Value = 1624.5
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H80").Value = 159885.6
$ws.Range("J80").Value = 159885.6
$ws.Range("L80").Value = 159885.6
$ws.Range("N80").Value = -161881.6
$ws.Range("H83").Value = 159885.6
$ws.Range("J83").Value = 159885.6
$ws.Range("L83").Value = 479656.8
$ws.Range("N83").Value = -489640.8
$ws.Range("H116").Value = 4330777.5
$ws.Range("I116").Value = 4786433
$ws.Range("J116").Value = 2050
$ws.Range("K116").Value = 4786433
$ws.Range("L116").Value = 2050
$ws.Range("M116").Value = -4784139
$ws.Range("N116").Value = -6638
$ws.Range("H122").Value = 4801.091
$ws.Range("I122").Value = 5118.1665
$ws.Range("J122").Value = 3374.25
$ws.Range("K122").Value = 15354.4995
$ws.Range("L122").Value = 10122.75
$ws.Range("M122").Value = -12904.4995
$ws.Range("N122").Value = -15022.75
$ws.Range("H132").Value = 5746.625
$ws.Range("I132").Value = 2911.3333
$ws.Range("K132").Value = 8733.999899999999
$ws.Range("M132").Value = -6203.999899999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4330777.5
$ws.Range("I3").Value = 4786433
$ws.Range("J3").Value = 2050
$ws.Range("K3").Value = 4786433
$ws.Range("L3").Value = 2050
$ws.Range("M3").Value = -4786319
$ws.Range("N3").Value = -2278
$ws.Range("H25").Value = 6222.7856
$ws.Range("I25").Value = 380.75
$ws.Range("J25").Value = 8559.6
$ws.Range("K25").Value = 380.75
$ws.Range("L25").Value = 8559.6
$ws.Range("M25").Value = -145.75
$ws.Range("N25").Value = -9029.6
$ws.Range("H94").Value = 646.86206
$ws.Range("I94").Value = 544.087
$ws.Range("J94").Value = 1040.8334
$ws.Range("K94").Value = 544.087
$ws.Range("L94").Value = 1040.8334
$ws.Range("M94").Value = -93.08699999999999
$ws.Range("N94").Value = -1942.8334
$ws.Range("H107").Value = 1734.1666
$ws.Range("J107").Value = 1962.5
$ws.Range("L107").Value = 1962.5
$ws.Range("N107").Value = -5802.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3016.7778
$ws.Range("I99").Value = 2999.2942
$ws.Range("J99").Value = 3314
$ws.Range("K99").Value = 2999.2942
$ws.Range("L99").Value = 3314
$ws.Range("M99").Value = -1501.2942
$ws.Range("N99").Value = -6310
$ws.Range("H107").Value = 2271.7273
$ws.Range("I107").Value = 1621.7778
$ws.Range("K107").Value = 1621.7778
$ws.Range("M107").Value = 298.2221999999999
$ws.Range("H122").Value = 1499.5
$ws.Range("I122").Value = 1499.5
$ws.Range("K122").Value = 4498.5
$ws.Range("M122").Value = -2048.5
$ws.Range("H126").Value = 3016.7778
$ws.Range("I126").Value = 2999.2942
$ws.Range("J126").Value = 3314
$ws.Range("K126").Value = 8997.882599999999
$ws.Range("L126").Value = 9942
$ws.Range("M126").Value = -6527.882599999999
$ws.Range("N126").Value = -14882
$ws.Range("H132").Value = 40736.324
$ws.Range("I132").Value = 1537.44
$ws.Range("K132").Value = 4612.32
$ws.Range("M132").Value = -2082.32

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 5855.909
$ws.Range("I18").Value = 500
$ws.Range("J18").Value = 20138.334
$ws.Range("K18").Value = 1500
$ws.Range("L18").Value = 60415.00199999999
$ws.Range("M18").Value = -1331
$ws.Range("N18").Value = -60753.00199999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H102").Value = 3532.3333
$ws.Range("I102").Value = 3532.3333
$ws.Range("K102").Value = 3532.3333
$ws.Range("M102").Value = -1910.3333
$ws.Range("H122").Value = 6009.7144
$ws.Range("I122").Value = 6009.7144
$ws.Range("K122").Value = 18029.1432
$ws.Range("M122").Value = -15579.1432
$ws.Range("H126").Value = 3285.9583
$ws.Range("I126").Value = 1940.0769
$ws.Range("J126").Value = 4876.5454
$ws.Range("K126").Value = 5820.2307
$ws.Range("L126").Value = 14629.6362
$ws.Range("M126").Value = -3350.2307
$ws.Range("N126").Value = -19569.6362
$ws.Range("H132").Value = 7740.514
$ws.Range("I132").Value = 6600.0586
$ws.Range("J132").Value = 8817.611000000001
$ws.Range("K132").Value = 19800.1758
$ws.Range("L132").Value = 26452.833
$ws.Range("M132").Value = -17270.1758
$ws.Range("N132").Value = -31512.833

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 11242.417
$ws.Range("I40").Value = 11091
$ws.Range("K40").Value = 11091
$ws.Range("M40").Value = -10955

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 39250
$ws.Range("J48").Value = 39250
$ws.Range("L48").Value = 39250
$ws.Range("N48").Value = -40388
$ws.Range("H49").Value = 67485.25
$ws.Range("J49").Value = 74961.664
$ws.Range("L49").Value = 74961.664
$ws.Range("N49").Value = -75421.664
$ws.Range("H107").Value = 1941.7428
$ws.Range("I107").Value = 2201.2593
$ws.Range("K107").Value = 6603.777900000001
$ws.Range("M107").Value = -4683.777900000001
$ws.Range("H122").Value = 3835.7896
$ws.Range("I122").Value = 3574.3572
$ws.Range("J122").Value = 4567.8
$ws.Range("K122").Value = 10723.0716
$ws.Range("L122").Value = 13703.4
$ws.Range("M122").Value = -8273.071599999999
$ws.Range("N122").Value = -18603.4
$ws.Range("H126").Value = 60126
$ws.Range("I126").Value = 63609.883
$ws.Range("J126").Value = 900
$ws.Range("K126").Value = 190829.649
$ws.Range("L126").Value = 2700
$ws.Range("M126").Value = -188359.649
$ws.Range("N126").Value = -7640
$ws.Range("H137").Value = 117899.664
$ws.Range("J137").Value = 117899.664
$ws.Range("L137").Value = 117899.664
$ws.Range("N137").Value = -128099.664
